$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1412156.9
$ws.Range("I15").Value = 1412156.9
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 4236470.699999999
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -4236301.699999999

# ALC row 17
$ws.Range("H17").Value = 1194.6964
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1194.6964
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3584.0892
$ws.Range("N17").Value = -3920.0892

# ALC row 32
$ws.Range("H32").Value = 1275
$ws.Range("I32").Value = 1275
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1275
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -949
$ws.Range("N32").ClearContents()

# ALC row 63
$ws.Range("H63").Value = 53950
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 53950
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 53950
$ws.Range("N63").Value = -55198

# ALC row 66
$ws.Range("H66").Value = 53950
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 53950
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 161850
$ws.Range("N66").Value = -168090

# ALC row 113
$ws.Range("H113").Value = 2658.7073
$ws.Range("I113").Value = 2111
$ws.Range("J113").Value = 3180.3333
$ws.Range("K113").Value = 2111
$ws.Range("L113").Value = 3180.3333
$ws.Range("M113").Value = 1143
$ws.Range("N113").Value = -9688.3333

# ALC row 129
$ws.Range("H129").Value = 1264.1708
$ws.Range("I129").Value = 508.875
$ws.Range("J129").Value = 1747.56
$ws.Range("K129").Value = 1526.625
$ws.Range("L129").Value = 5242.68
$ws.Range("M129").Value = 3473.375
$ws.Range("N129").Value = -15242.68

# ALC row 132
$ws.Range("H132").Value = 8681.079
$ws.Range("I132").Value = 10689.723
$ws.Range("J132").Value = 6873.3
$ws.Range("K132").Value = 32069.169
$ws.Range("L132").Value = 20619.9
$ws.Range("M132").Value = -29539.169
$ws.Range("N132").Value = -25679.9

# ALC row 137
$ws.Range("H137").Value = 1278.5217
$ws.Range("I137").Value = 1621.8
$ws.Range("J137").Value = 925.14703
$ws.Range("K137").Value = 4865.4
$ws.Range("L137").Value = 2775.44109
$ws.Range("M137").Value = -2315.4
$ws.Range("N137").Value = -7875.44109

# ALC row 138
$ws.Range("H138").Value = 1273.9828
$ws.Range("I138").Value = 1097.2812
$ws.Range("J138").Value = 1491.4615
$ws.Range("K138").Value = 3291.8436
$ws.Range("L138").Value = 4474.3845
$ws.Range("M138").Value = 1848.1564
$ws.Range("N138").Value = -14754.3845

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3502.52
$ws.Range("I32").Value = 3226.5684
$ws.Range("J32").Value = 8745.6
$ws.Range("K32").Value = 3226.5684
$ws.Range("L32").Value = 8745.6
$ws.Range("M32").Value = -2939.5684
$ws.Range("N32").Value = -9319.6

# ARM row 57
$ws.Range("H57").Value = 14375
$ws.Range("I57").Value = 14375
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 14375
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -13891

# ARM row 122
$ws.Range("H122").Value = 1343
$ws.Range("I122").Value = 1065
$ws.Range("J122").Value = 2246.5
$ws.Range("K122").Value = 3195
$ws.Range("L122").Value = 6739.5
$ws.Range("M122").Value = -745
$ws.Range("N122").Value = -11639.5

# BSM row 57
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 39600
$ws.Range("I57").Value = 10000
$ws.Range("J57").Value = 54400
$ws.Range("K57").Value = 10000
$ws.Range("L57").Value = 54400
$ws.Range("M57").Value = -9280
$ws.Range("N57").Value = -55840

# BSM row 105
$ws.Range("H105").Value = 1502
$ws.Range("I105").Value = 1323.3529
$ws.Range("J105").Value = 1718.9286
$ws.Range("K105").Value = 1323.3529
$ws.Range("L105").Value = 1718.9286
$ws.Range("M105").Value = 423.6470999999999
$ws.Range("N105").Value = -5212.9286

# BSM row 136
$ws.Range("H136").Value = 39600
$ws.Range("I136").Value = 10000
$ws.Range("J136").Value = 54400
$ws.Range("K136").Value = 10000
$ws.Range("L136").Value = 54400
$ws.Range("M136").Value = -4900
$ws.Range("N136").Value = -64600

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4833016.5
$ws.Range("I31").Value = 1581.5103
$ws.Range("J31").Value = 16670033
$ws.Range("K31").Value = 1581.5103
$ws.Range("L31").Value = 16670033
$ws.Range("M31").Value = -1286.5103
$ws.Range("N31").Value = -16670623

# CRP row 34
$ws.Range("H34").Value = 4833016.5
$ws.Range("I34").Value = 1581.5103
$ws.Range("J34").Value = 16670033
$ws.Range("K34").Value = 1581.5103
$ws.Range("L34").Value = 16670033
$ws.Range("M34").Value = -1379.5103
$ws.Range("N34").Value = -16670437

# CRP row 58
$ws.Range("H58").Value = 772184.4
$ws.Range("I58").Value = 1567.6757
$ws.Range("J58").Value = 1790499.4
$ws.Range("K58").Value = 1567.6757
$ws.Range("L58").Value = 1790499.4
$ws.Range("M58").Value = -1364.6757
$ws.Range("N58").Value = -1790905.4

# CRP row 136
$ws.Range("H136").Value = 772184.4
$ws.Range("I136").Value = 1567.6757
$ws.Range("J136").Value = 1790499.4
$ws.Range("K136").Value = 4703.0271
$ws.Range("L136").Value = 5371498.199999999
$ws.Range("M136").Value = -2153.0271
$ws.Range("N136").Value = -5376598.199999999

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 707.55554
$ws.Range("I5").Value = 304.30768
$ws.Range("J5").Value = 1756
$ws.Range("K5").Value = 912.92304
$ws.Range("L5").Value = 5268
$ws.Range("M5").Value = -800.92304
$ws.Range("N5").Value = -5492

# CUL row 40
$ws.Range("H40").Value = 91.75
$ws.Range("I40").Value = 90.28570999999999
$ws.Range("J40").Value = 102
$ws.Range("K40").Value = 361.14284
$ws.Range("L40").Value = 408
$ws.Range("M40").Value = -292.14284
$ws.Range("N40").Value = -546

# CUL row 46
$ws.Range("H46").Value = 691.0909
$ws.Range("I46").Value = 622.44446
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 1867.33338
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -1776.33338
$ws.Range("N46").Value = -3182

# CUL row 57
$ws.Range("H57").Value = 3833.3333
$ws.Range("I57").Value = 3833.3333
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 11499.9999
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -10940.9999

# CUL row 58
$ws.Range("H58").Value = 1723.75
$ws.Range("I58").Value = 965
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 2895
$ws.Range("L58").Value = 12000
$ws.Range("M58").Value = -2767
$ws.Range("N58").Value = -12256

# CUL row 122
$ws.Range("H122").Value = 2753.4849
$ws.Range("I122").Value = 520.2
$ws.Range("J122").Value = 3152.2856
$ws.Range("K122").Value = 4681.8
$ws.Range("L122").Value = 28370.5704
$ws.Range("M122").Value = -2231.8
$ws.Range("N122").Value = -33270.5704

# CUL row 135
$ws.Range("H135").Value = 707.55554
$ws.Range("I135").Value = 304.30768
$ws.Range("J135").Value = 1756
$ws.Range("K135").Value = 2738.76912
$ws.Range("L135").Value = 15804
$ws.Range("M135").Value = -203.7691199999999
$ws.Range("N135").Value = -20874

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5459.3335
$ws.Range("I70").Value = 4600.6875
$ws.Range("J70").Value = 6440.643
$ws.Range("K70").Value = 4600.6875
$ws.Range("L70").Value = 6440.643
$ws.Range("M70").Value = -4330.6875
$ws.Range("N70").Value = -6980.643

# GSM row 73
$ws.Range("H73").Value = 5459.3335
$ws.Range("I73").Value = 4600.6875
$ws.Range("J73").Value = 6440.643
$ws.Range("K73").Value = 4600.6875
$ws.Range("L73").Value = 6440.643
$ws.Range("M73").Value = -3664.6875
$ws.Range("N73").Value = -8312.643

# GSM row 80
$ws.Range("H80").Value = 3141.5789
$ws.Range("I80").Value = 2587.2222
$ws.Range("J80").Value = 3640.5
$ws.Range("K80").Value = 2587.2222
$ws.Range("L80").Value = 3640.5
$ws.Range("M80").Value = -1589.2222
$ws.Range("N80").Value = -5636.5

# GSM row 83
$ws.Range("H83").Value = 3141.5789
$ws.Range("I83").Value = 2587.2222
$ws.Range("J83").Value = 3640.5
$ws.Range("K83").Value = 12936.111
$ws.Range("L83").Value = 18202.5
$ws.Range("M83").Value = -7944.111000000001
$ws.Range("N83").Value = -28186.5

# GSM row 102
$ws.Range("H102").Value = 1849.6875
$ws.Range("I102").Value = 1644
$ws.Range("J102").Value = 2242.3635
$ws.Range("K102").Value = 1644
$ws.Range("L102").Value = 2242.3635
$ws.Range("M102").Value = -22
$ws.Range("N102").Value = -5486.363499999999

# GSM row 107
$ws.Range("H107").Value = 737.16
$ws.Range("I107").Value = 640.2857
$ws.Range("J107").Value = 1245.75
$ws.Range("K107").Value = 640.2857
$ws.Range("L107").Value = 1245.75
$ws.Range("M107").Value = 1279.7143
$ws.Range("N107").Value = -5085.75

# LTW row 9
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 574.2857
$ws.Range("I9").Value = 170
$ws.Range("J9").Value = 3000
$ws.Range("K9").Value = 170
$ws.Range("L9").Value = 3000
$ws.Range("M9").Value = 54
$ws.Range("N9").Value = -3448

# LTW row 132
$ws.Range("H132").Value = 15281.16
$ws.Range("I132").Value = 22830.213
$ws.Range("J132").Value = 2609.5356
$ws.Range("K132").Value = 68490.639
$ws.Range("L132").Value = 7828.6068
$ws.Range("M132").Value = -65960.639
$ws.Range("N132").Value = -12888.6068

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2406.8276
$ws.Range("I126").Value = 2595.75
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 7787.25
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = -5317.25
$ws.Range("N126").Value = -9440

# WVR row 132
$ws.Range("H132").Value = 2213.92
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2213.92
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 6641.76
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -11701.76
